$wb = $excel.ActiveWorkbook

# Update "展览" sheet: F7 6652 -> 6653, F16 213 -> 214
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 6653
$ws1.Range("F16").Value = 214

# Update "全部类型" sheet: F7 6652 -> 6653, F16 213 -> 214
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 6653
$ws4.Range("F16").Value = 214
